# Weekly refresh of the "Fruta, Vega Modelo de Temuco - Pera asiática" price
# sheet: the data rows (2-18) get reshuffled onto new dates/prices — each
# row ends up carrying the Fecha/Calidad/Volumen/Precio*/Unidad/Origen/
# Precio-Kg/Kg-unidad values that (before this edit) belonged to another
# row in the same block. Row 19 and the header row are untouched, as are
# the constant identifying columns (A-C, E-K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" for each row.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D,L,M,N,O,P,Q,R,S,T

# Snapshot the current (pre-edit) values for every data row before writing
# anything, so later writes never read already-mutated data.
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destination row -> source row (where its new D/L/M/N/O/P/Q/R/S/T values
# come from, read out of the snapshot taken above).
$mapping = @{
    2  = 4
    3  = 10
    4  = 13
    5  = 12
    6  = 8
    7  = 17
    8  = 5
    9  = 16
    10 = 3
    11 = 15
    12 = 6
    13 = 7
    14 = 9
    15 = 14
    16 = 11
    17 = 18
    18 = 2
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
